# Client update for experience page
# - Update the Banyan Point Condominiums description text (row 41) and
#   grow its row height to fit the extra wrapped line.
# - Add a new client row (row 42) for "Upwork & Fiverr" / Freelance work,
#   including a real date value in the client_start_date column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41: Banyan Point Condominiums - refreshed description ---
$ws.Range("E41").Value = "A Squarespace site for a condominium HOA. Assisted in rebuilding links, creating/editing and uploading assets to their Squarespace site. Providing general technical support with PDF and Microsoft Office products."
$ws.Rows.Item(41).RowHeight = 43.2

# --- Row 42: new Freelance client entry ---
$ws.Range("A42").Value = 41
$ws.Range("B42").Value = "Upwork & Fiverr"
$ws.Range("C42").Value = "no"
$ws.Range("D42").Value = "Freelance"
$ws.Range("E42").Value = "Freelance work from Marketing Automation integration support to Full-stack Code Reviews and optimizations and Security Audits."
$ws.Range("E42").WrapText = $true
$ws.Range("F42").Value = "Web Development, HTML, JavaScript, jQuery, CSS, VB.Net, C#.Net, ASP.Net, Entity Framework, LINQ, SQL, T-SQL, MySQL, Oracle, Jira, AWS, React, Angular, JSP,  VisualForce Pages, ColdFusion, Git/GitHub, PHP, Java Servlets, RESTful Services"
$ws.Range("F42").WrapText = $true
$ws.Range("H42").Value = (Get-Date -Year 2025 -Month 10 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("H42").NumberFormat = "m/d/yyyy"
$ws.Range("J42").Value = "Freelance Work"
$ws.Rows.Item(42).RowHeight = 57.6

# Keep the view parked where the new rows were edited, like the live session.
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("F37").Select() | Out-Null
